# Weekly update: insert a new daily price record for "Choclo" (Choclero,
# Primera, Región del Maule) on 2022-02-01 (serial date 44606) at row 116.
# This pushes all the existing records that were at rows 116-190 down to
# rows 117-191, growing the used range from A1:R190 to A1:R191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116, shifting rows 116:190 down to 117:191.
$ws.Rows("116:116").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(116, 1).Value = 5
$ws.Cells.Item(116, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(116, 3).Value = "Maule"
$ws.Cells.Item(116, 4).Value = 44606
$ws.Cells.Item(116, 5).Value = 7
$ws.Cells.Item(116, 6).Value = 100112024
$ws.Cells.Item(116, 7).Value = "Choclo"
$ws.Cells.Item(116, 8).Value = "Choclero"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 40000
$ws.Cells.Item(116, 11).Value = 120
$ws.Cells.Item(116, 12).Value = 120
$ws.Cells.Item(116, 13).Value = 120
$ws.Cells.Item(116, 14).Value = "`$/unidad"
$ws.Cells.Item(116, 15).Value = "Región del Maule"
$ws.Cells.Item(116, 16).Value = 120
$ws.Cells.Item(116, 17).Value = 1
$ws.Cells.Item(116, 18).Value = "Hortaliza"
